$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to match the new workbook title.
$ws.Name = "Daily Progress Report"

# Keep the print title rows in sync with the sheet's new name
# (the defined name _xlnm.Print_Titles is sheet-qualified).
$ws.PageSetup.PrintTitleRows = "`$3:`$3"

# Reflect the last active cell selection recorded when the workbook was saved.
$ws.Range("D3").Select()
